# Applies the corrections described by the diff:
#  - Rows 58 and 59 had their match data (columns F:V) swapped
#  - Rows 63, 64 and 65 had their match data (columns F:V) rotated
#  - Rows 68 and 69 had their match data (columns F:V) swapped
#  - Three new match rows (117, 118, 119) were appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 58 <-> 59 (columns F:V) ---------------------------------
$v58 = $ws.Range("F58:V58").Value2()
$v59 = $ws.Range("F59:V59").Value2()
$ws.Range("F58:V58").Value2 = $v59
$ws.Range("F59:V59").Value2 = $v58

# --- Rotate rows 63 -> 64 -> 65 -> 63 (columns F:V) ---------------------
$v63 = $ws.Range("F63:V63").Value2()
$v64 = $ws.Range("F64:V64").Value2()
$v65 = $ws.Range("F65:V65").Value2()
$ws.Range("F63:V63").Value2 = $v64
$ws.Range("F64:V64").Value2 = $v65
$ws.Range("F65:V65").Value2 = $v63

# --- Swap rows 68 <-> 69 (columns F:V) ----------------------------------
$v68 = $ws.Range("F68:V68").Value2()
$v69 = $ws.Range("F69:V69").Value2()
$ws.Range("F68:V68").Value2 = $v69
$ws.Range("F69:V69").Value2 = $v68

# --- Append three new match rows (117-119), matching row 116's style ---
$ws.Range("A116:V116").Copy()
$ws.Range("A117:V119").PasteSpecial(-4122)

$ws.Range("A117").Value2 = 116
$ws.Range("B117").Value2 = "spain"
$ws.Range("C117").Value2 = "laliga"
$ws.Range("D117").Value2 = "2023-2024"
$ws.Range("E117").Value2 = 45235.67708333334
$ws.Range("F117").Value2 = "Valencia"
$ws.Range("G117").Value2 = 1
$ws.Range("H117").Value2 = "Granada CF"
$ws.Range("I117").Value2 = 0
$ws.Range("J117").Value2 = 1.54
$ws.Range("K117").Value2 = "22/10/2023 12:02"
$ws.Range("L117").Value2 = 1.67
$ws.Range("M117").Value2 = "05/11/2023 16:13"
$ws.Range("N117").Value2 = 4.19
$ws.Range("O117").Value2 = "22/10/2023 12:02"
$ws.Range("P117").Value2 = 4.15
$ws.Range("Q117").Value2 = "05/11/2023 16:13"
$ws.Range("R117").Value2 = 5.77
$ws.Range("S117").Value2 = "22/10/2023 12:02"
$ws.Range("T117").Value2 = 5.29
$ws.Range("U117").Value2 = "05/11/2023 16:13"
$ws.Range("V117").Value2 = "https://www.betexplorer.com/football/spain/laliga/valencia-granada-cf/MaSXAFrI/"

$ws.Range("A118").Value2 = 117
$ws.Range("B118").Value2 = "spain"
$ws.Range("C118").Value2 = "laliga"
$ws.Range("D118").Value2 = "2023-2024"
$ws.Range("E118").Value2 = 45235.77083333334
$ws.Range("F118").Value2 = "Villarreal"
$ws.Range("G118").Value2 = 2
$ws.Range("H118").Value2 = "Ath Bilbao"
$ws.Range("I118").Value2 = 3
$ws.Range("J118").Value2 = 2.5
$ws.Range("K118").Value2 = "22/10/2023 12:02"
$ws.Range("L118").Value2 = 3.82
$ws.Range("M118").Value2 = "05/11/2023 18:25"
$ws.Range("N118").Value2 = 3.56
$ws.Range("O118").Value2 = "22/10/2023 12:02"
$ws.Range("P118").Value2 = 3.94
$ws.Range("Q118").Value2 = "05/11/2023 18:29"
$ws.Range("R118").Value2 = 2.68
$ws.Range("S118").Value2 = "22/10/2023 12:02"
$ws.Range("T118").Value2 = 1.95
$ws.Range("U118").Value2 = "05/11/2023 18:29"
$ws.Range("V118").Value2 = "https://www.betexplorer.com/football/spain/laliga/villarreal-ath-bilbao/n75GGJEt/"

$ws.Range("A119").Value2 = 118
$ws.Range("B119").Value2 = "spain"
$ws.Range("C119").Value2 = "laliga"
$ws.Range("D119").Value2 = "2023-2024"
$ws.Range("E119").Value2 = 45235.875
$ws.Range("F119").Value2 = "Real Madrid"
$ws.Range("G119").Value2 = 0
$ws.Range("H119").Value2 = "Rayo Vallecano"
$ws.Range("I119").Value2 = 0
$ws.Range("J119").Value2 = 1.3
$ws.Range("K119").Value2 = "22/10/2023 12:02"
$ws.Range("L119").Value2 = 1.22
$ws.Range("M119").Value2 = "05/11/2023 20:38"
$ws.Range("N119").Value2 = 5.81
$ws.Range("O119").Value2 = "22/10/2023 12:02"
$ws.Range("P119").Value2 = 6.9
$ws.Range("Q119").Value2 = "05/11/2023 20:59"
$ws.Range("R119").Value2 = 10.18
$ws.Range("S119").Value2 = "22/10/2023 12:02"
$ws.Range("T119").Value2 = 13.8
$ws.Range("U119").Value2 = "05/11/2023 20:59"
$ws.Range("V119").Value2 = "https://www.betexplorer.com/football/spain/laliga/real-madrid-rayo-vallecano/r9TTBeTB/"

Write-Output "edit applied"
